$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Maracuyá (Agrícola del Norte S.A. de
# Arica) on 2022-01-14 (Excel serial 44575). Insert it as a new row 18,
# pushing the existing rows 18-76 down to 19-77.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44575
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100108
$ws.Range("H18").Value = "Tropicales y subtropicales"
$ws.Range("I18").Value = 100108003
$ws.Range("J18").Value = "Maracuyá"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 54000
$ws.Range("O18").Value = 55000
$ws.Range("P18").Value = 54500
$ws.Range("Q18").Value = "$/caja 20 kilos"
$ws.Range("R18").Value = "Región de Arica y Parinacota"
$ws.Range("S18").Value = 2725
$ws.Range("T18").Value = 20
